$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="31÷5=6, 1"},
    @{Row=1;  Col=2; New="85÷6=14, 1"},
    @{Row=1;  Col=3; New="69÷3=23, 0"},
    @{Row=1;  Col=4; New="44÷4=11, 0"},
    @{Row=1;  Col=5; New="73÷3=24, 1"},

    @{Row=5;  Col=1; New="91÷6=15, 1"},
    @{Row=5;  Col=2; New="52÷4=13, 0"},
    @{Row=5;  Col=3; New="14÷4=3, 2"},
    @{Row=5;  Col=4; New="19÷3=6, 1"},
    @{Row=5;  Col=5; New="74÷3=24, 2"},

    @{Row=9;  Col=1; New="35÷5=7, 0"},
    @{Row=9;  Col=2; New="67÷8=8, 3"},
    @{Row=9;  Col=3; New="45÷5=9, 0"},
    @{Row=9;  Col=4; New="32÷3=10, 2"},
    @{Row=9;  Col=5; New="23÷9=2, 5"},

    @{Row=13; Col=1; New="54÷4=13, 2"},
    @{Row=13; Col=2; New="29÷2=14, 1"},
    @{Row=13; Col=3; New="23÷6=3, 5"},
    @{Row=13; Col=4; New="80÷7=11, 3"},
    @{Row=13; Col=5; New="57÷4=14, 1"},

    @{Row=17; Col=1; New="79÷9=8, 7"},
    @{Row=17; Col=2; New="76÷6=12, 4"},
    @{Row=17; Col=3; New="54÷9=6, 0"},
    @{Row=17; Col=4; New="12÷8=1, 4"},
    @{Row=17; Col=5; New="54÷9=6, 0"}
)

foreach ($u in $updates) {
    $cell = $tbl.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $u.New
}
